# Generate Report for Handoff
# Adds a new handed-off file (5a9ca048-9834-4a26-9179-dfd3d34d708b.md) as a
# new row on each of the three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$commit = "5d89a0a8a596dab22ff620aceef655f57159b218"
$newFile = "5a9ca048-9834-4a26-9179-dfd3d34d708b.md"
$newFileUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/$commit/e2e/$newFile"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Columns.Item(1).ColumnWidth = 39.14

$wsOverview.Cells.Item(3, 1).Value = $newFile
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 4).Value = ""
$wsOverview.Cells.Item(3, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 7).Value = "2016-08-13 08:46:49"

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3, 2), $newFileUrl, "", "", "e2e\$newFile") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Columns.Item(1).ColumnWidth = 39.14

$wsZhCn.Cells.Item(3, 2).Value = ".md"
$wsZhCn.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(3, 4).Value = "e2e"
$wsZhCn.Cells.Item(3, 5).Value = "ht"
$wsZhCn.Cells.Item(3, 6).Value = "False"
$wsZhCn.Cells.Item(3, 7).Value = "5a9ca048-9834-4a26-9179-dfd3d34d708b.e75ae6a34a82b7cb27aab8e848fdcd1d0fccaaac.zh-cn.xlf"
$wsZhCn.Cells.Item(3, 8).Value = "2016-08-13 08:46:41"
$wsZhCn.Cells.Item(3, 9).Value = ""
$wsZhCn.Cells.Item(3, 10).Value = ""
$wsZhCn.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(3, 12).Value = ""
$wsZhCn.Cells.Item(3, 13).Value = "True"
$wsZhCn.Cells.Item(3, 14).Value = ""
$wsZhCn.Cells.Item(3, 15).Value = "False"
$wsZhCn.Cells.Item(3, 16).Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(3, 1), $newFileUrl, "", "", $newFile) | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Columns.Item(1).ColumnWidth = 39.14

$wsDeDe.Cells.Item(3, 2).Value = ".md"
$wsDeDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(3, 4).Value = "e2e"
$wsDeDe.Cells.Item(3, 5).Value = "ht"
$wsDeDe.Cells.Item(3, 6).Value = "False"
$wsDeDe.Cells.Item(3, 7).Value = "5a9ca048-9834-4a26-9179-dfd3d34d708b.e75ae6a34a82b7cb27aab8e848fdcd1d0fccaaac.de-de.xlf"
$wsDeDe.Cells.Item(3, 8).Value = "2016-08-13 08:46:49"
$wsDeDe.Cells.Item(3, 9).Value = ""
$wsDeDe.Cells.Item(3, 10).Value = ""
$wsDeDe.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(3, 12).Value = ""
$wsDeDe.Cells.Item(3, 13).Value = "True"
$wsDeDe.Cells.Item(3, 14).Value = ""
$wsDeDe.Cells.Item(3, 15).Value = "False"
$wsDeDe.Cells.Item(3, 16).Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(3, 1), $newFileUrl, "", "", $newFile) | Out-Null
